# The document has one section whose first-page header/footer and
# default (odd-page) header/footer each carry a logo picture as an
# inline picture (wp:inline / pic:pic). The edit simply relabels the
# picture's stored "name" (wp:docPr / pic:cNvPr @name) on three of
# these inline pictures - no visual / layout change.
#
# Mapping (confirmed via word/_rels/document.xml.rels +
# w:headerReference/w:footerReference types):
#   header1.xml -> Sections(1).Headers(wdHeaderFooterFirstPage)   BTec logo   image1.jpg -> image2.jpg
#   footer1.xml -> Sections(1).Footers(wdHeaderFooterFirstPage)   Pearson logo id=3  image2.png -> image1.png
#   footer2.xml -> Sections(1).Footers(wdHeaderFooterPrimary)     Pearson logo id=2  image2.png -> image1.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3

# footer2.xml (default/primary footer) - Pearson logo, id=2
$ftrPrimary = $sec.Footers.Item(1)
$ftrPrimary.Range.InlineShapes.Item(1).Name = "image1.png"

# header1.xml (first-page header) - BTec logo, id=1
$hdrFirst = $sec.Headers.Item(2)
$hdrFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"

# footer1.xml (first-page footer) - Pearson logo, id=3
$ftrFirst = $sec.Footers.Item(2)
$ftrFirst.Range.InlineShapes.Item(1).Name = "image1.png"
